$d = $word.ActiveDocument
$successCount = 0
$failCount = 0

# Replacement 1
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("enbolbí mayornan aktivamente", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $rng1.Select()
    $word.Selection.TypeText("enbolbí mayornan/edukadónan aktivamente")
    $successCount++
} else {
    Write-Host "FAILED to find replacement 1: enbolbí mayornan aktivamente"
    $failCount++
}

# Replacement 2
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("fasilitá un diskushon ku mayornan:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $rng2.Select()
    $word.Selection.TypeText("fasilitá un diskushon ku mayornan/edukadónan:")
    $successCount++
} else {
    Write-Host "FAILED to find replacement 2: fasilitá un diskushon ku mayornan:"
    $failCount++
}

# Replacement 3
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("pa asina mayornan kuminsá hasi", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $rng3.Select()
    $word.Selection.TypeText("pa asina mayornan/edukadónan kuminsá hasi")
    $successCount++
} else {
    Write-Host "FAILED to find replacement 3: pa asina mayornan kuminsá hasi"
    $failCount++
}

# Replacement 4
$rng4 = $d.Content
$found4 = $rng4.Find.Execute("Aseptá kontribushon di mayornan ta enkurashá", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found4) {
    $rng4.Select()
    $word.Selection.TypeText("Aseptá kontribushon di mayornan/edukadónan ta enkurashá")
    $successCount++
} else {
    Write-Host "FAILED to find replacement 4: Aseptá kontribushon di mayornan ta enkur"
    $failCount++
}

# Replacement 5
$rng5 = $d.Content
$found5 = $rng5.Find.Execute("klave pa mayornan duna ehèmpel", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found5) {
    $rng5.Select()
    $word.Selection.TypeText("klave pa mayornan/edukadónan duna ehèmpel")
    $successCount++
} else {
    Write-Host "FAILED to find replacement 5: klave pa mayornan duna ehèmpel"
    $failCount++
}

# Replacement 6
$rng6 = $d.Content
$found6 = $rng6.Find.Execute("loke e mayor ta bisa pa mustra", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found6) {
    $rng6.Select()
    $word.Selection.TypeText("loke e mayor/edukadó ta bisa pa mustra")
    $successCount++
} else {
    Write-Host "FAILED to find replacement 6: loke e mayor ta bisa pa mustra"
    $failCount++
}

# Replacement 7
$rng7 = $d.Content
$found7 = $rng7.Find.Execute("pa ta e mayor ku ta permití bo", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found7) {
    $rng7.Select()
    $word.Selection.TypeText("pa ta e mayor/edukadó ku ta permití bo")
    $successCount++
} else {
    Write-Host "FAILED to find replacement 7: pa ta e mayor ku ta permití bo"
    $failCount++
}

# Replacement 8
$rng8 = $d.Content
$found8 = $rng8.Find.Execute("Mayor: Mi no tabata sintí", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found8) {
    $rng8.Select()
    $word.Selection.TypeText("Mayor/Edukadó: Mi no tabata sintí")
    $successCount++
} else {
    Write-Host "FAILED to find replacement 8: Mayor: Mi no tabata sintí"
    $failCount++
}

# Replacement 9
$rng9 = $d.Content
$found9 = $rng9.Find.Execute("eksperensia di un mayor den mas detaye", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found9) {
    $rng9.Select()
    $word.Selection.TypeText("eksperensia di un mayor/edukadó den mas detaye")
    $successCount++
} else {
    Write-Host "FAILED to find replacement 9: eksperensia di un mayor den mas detaye"
    $failCount++
}

# Replacement 10
$rng10 = $d.Content
$found10 = $rng10.Find.Execute("yuda mayornan reflehá riba kon", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found10) {
    $rng10.Select()
    $word.Selection.TypeText("yuda mayornan/edukadónan reflehá riba kon")
    $successCount++
} else {
    Write-Host "FAILED to find replacement 10: yuda mayornan reflehá riba kon"
    $failCount++
}

# Replacement 11
$rng11 = $d.Content
$found11 = $rng11.Find.Execute("emoshonnan di e mayor komo di e yu", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found11) {
    $rng11.Select()
    $word.Selection.TypeText("emoshonnan di e mayor/edukadó komo di e yu")
    $successCount++
} else {
    Write-Host "FAILED to find replacement 11: emoshonnan di e mayor komo di e yu"
    $failCount++
}

# Replacement 12
$rng12 = $d.Content
$found12 = $rng12.Find.Execute("Kon sigui bo yu su guia a lagá bo sinti? Kon bo ta kere ku el a laga bo yu sinti? ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found12) {
    $rng12.Select()
    $word.Selection.TypeText("Kon bo a sinti pa sigui e yu su guia? Kon bo ta kere ku bo yu a sinti? ")
    $successCount++
} else {
    Write-Host "FAILED to find replacement 12: Kon sigui bo yu su guia a lagá bo sinti?"
    $failCount++
}

# Replacement 13
$rng13 = $d.Content
$found13 = $rng13.Find.Execute("Konektá e eksperensia di e mayor ku prinsipionan", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found13) {
    $rng13.Select()
    $word.Selection.TypeText("Konektá e eksperensia di e mayor/edukadó ku prinsipionan")
    $successCount++
} else {
    Write-Host "FAILED to find replacement 13: Konektá e eksperensia di e mayor ku prin"
    $failCount++
}

# Replacement 14
$rng14 = $d.Content
$found14 = $rng14.Find.Execute("Esaki ta yuda mayornan komprondé nan eksperensia", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found14) {
    $rng14.Select()
    $word.Selection.TypeText("Esaki ta yuda mayornan/edukadónan komprondé nan eksperensia")
    $successCount++
} else {
    Write-Host "FAILED to find replacement 14: Esaki ta yuda mayornan komprondé nan eks"
    $failCount++
}

# Replacement 15
$rng15 = $d.Content
$found15 = $rng15.Find.Execute("Enkurashá mayornan pa krea nan mes konekshonnan", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found15) {
    $rng15.Select()
    $word.Selection.TypeText("Enkurashá mayornan/edukadónan pa krea nan mes konekshonnan")
    $successCount++
} else {
    Write-Host "FAILED to find replacement 15: Enkurashá mayornan pa krea nan mes konek"
    $failCount++
}

# Replacement 16
$rng16 = $d.Content
$found16 = $rng16.Find.Execute("Para ketu i pone énfasis riba patronchinan òf temanan den vários historia di mayor", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found16) {
    $rng16.Select()
    $word.Selection.TypeText("Enfatisá ('Highlight') patronchinan òf temanan den e diferente historianan di mayor/edukadó ")
    $successCount++
} else {
    Write-Host "FAILED to find replacement 16: Para ketu i pone énfasis riba patronchin"
    $failCount++
}

# Replacement 17
$rng17 = $d.Content
$found17 = $rng17.Find.Execute("Referí bèk na diskushonnan òf plachinan anterior", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found17) {
    $rng17.Select()
    $word.Selection.TypeText("Referí na diskushonnan òf plachinan anterior")
    $successCount++
} else {
    Write-Host "FAILED to find replacement 17: Referí bèk na diskushonnan òf plachinan "
    $failCount++
}

# Replacement 18
$rng18 = $d.Content
$found18 = $rng18.Find.Execute("Esaki ta bai bèk na loke nos a siña", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found18) {
    $rng18.Select()
    $word.Selection.TypeText("Esaki ta hiba nos bèk na loke nos a siña")
    $successCount++
} else {
    Write-Host "FAILED to find replacement 18: Esaki ta bai bèk na loke nos a siña"
    $failCount++
}

# Replacement 19
$rng19 = $d.Content
$found19 = $rng19.Find.Execute("Evidensia sientífiko ta bisa nos ku mayornan tin mas chèns", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found19) {
    $rng19.Select()
    $word.Selection.TypeText("Evidensia sientífiko ta bisa nos ku mayornan/edukadónan tin mas chèns")
    $successCount++
} else {
    Write-Host "FAILED to find replacement 19: Evidensia sientífiko ta bisa nos ku mayo"
    $failCount++
}

# Replacement 20
$rng20 = $d.Content
$found20 = $rng20.Find.Execute("bo mester duna mayornan mas tantu oportunidat", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found20) {
    $rng20.Select()
    $word.Selection.TypeText("bo mester duna mayornan/edukadónan mas tantu oportunidat")
    $successCount++
} else {
    Write-Host "FAILED to find replacement 20: bo mester duna mayornan mas tantu oportu"
    $failCount++
}

# Replacement 21
$rng21 = $d.Content
$found21 = $rng21.Find.Execute("manera ora un mayor ta kompartí un situashon difísil. Esaki ta permití e mayor pa praktiká", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found21) {
    $rng21.Select()
    $word.Selection.TypeText("manera ora un mayor/edukadó ta kompartí un situashon difísil. Esaki ta permití e mayor/edukadó pa praktiká")
    $successCount++
} else {
    Write-Host "FAILED to find replacement 21: manera ora un mayor ta kompartí un situa"
    $failCount++
}

Write-Host "Done: $successCount succeeded, $failCount failed"
